$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1871727748691099
$ws.Range("C2").Value = 0.569371727748691
$ws.Range("J2").Value = 0.02486910994764398
$ws.Range("O2").Value = 0.00130890052356021
$ws.Range("P2").Value = 0.1282722513089005
$ws.Range("S2").Value = 0.08900523560209424
$ws.Range("B3").Value = 0.01138952164009112
$ws.Range("C3").Value = 0.02050113895216401
$ws.Range("J3").Value = 0.02733485193621868
$ws.Range("P3").Value = 0.744874715261959
$ws.Range("S3").Value = 0.1958997722095672
$ws.Range("J4").Value = 0.05660377358490566
$ws.Range("P4").Value = 0.6698113207547169
$ws.Range("S4").Value = 0.2735849056603774
$ws.Range("B6").Value = 0.1027397260273973
$ws.Range("D6").Value = 0.02054794520547945
$ws.Range("F6").Value = 0.0593607305936073
$ws.Range("J6").Value = 0.2488584474885845
$ws.Range("O6").Value = 0.0182648401826484
$ws.Range("Q6").Value = 0.1438356164383562
$ws.Range("R6").Value = 0.0547945205479452
$ws.Range("S6").Value = 0.3515981735159817
$ws.Range("B7").Value = 0.1396648044692737
$ws.Range("D7").Value = 0.0260707635009311
$ws.Range("E7").Value = 0.0037243947858473
$ws.Range("F7").Value = 0.0446927374301676
$ws.Range("J7").Value = 0.1284916201117318
$ws.Range("O7").Value = 0.01303538175046555
$ws.Range("Q7").Value = 0.1918063314711359
$ws.Range("R7").Value = 0.08193668528864059
$ws.Range("S7").Value = 0.3705772811918063
$ws.Range("B8").Value = 0.1202651515151515
$ws.Range("D8").Value = 0.01988636363636364
$ws.Range("F8").Value = 0.04450757575757576
$ws.Range("J8").Value = 0.1174242424242424
$ws.Range("O8").Value = 0.01041666666666667
$ws.Range("Q8").Value = 0.15625
$ws.Range("R8").Value = 0.09659090909090909
$ws.Range("S8").Value = 0.4346590909090909
$ws.Range("B9").Value = 0.1114206128133705
$ws.Range("D9").Value = 0.02785515320334262
$ws.Range("F9").Value = 0.07799442896935933
$ws.Range("J9").Value = 0.116991643454039
$ws.Range("O9").Value = 0.01114206128133705
$ws.Range("Q9").Value = 0.1838440111420613
$ws.Range("R9").Value = 0.1002785515320334
$ws.Range("S9").Value = 0.3704735376044568
$ws.Range("B10").Value = 0.1248550444530344
$ws.Range("D10").Value = 0.02319288751449556
$ws.Range("E10").Value = 0.00154619250096637
$ws.Range("F10").Value = 0.06725937379203711
$ws.Range("J10").Value = 0.1140316969462698
$ws.Range("O10").Value = 0.01623502126014689
$ws.Range("Q10").Value = 0.2060301507537688
$ws.Range("R10").Value = 0.09006571318129107
$ws.Range("S10").Value = 0.3567839195979899
$ws.Range("G11").Value = 0.1492957746478873
$ws.Range("J11").Value = 0.08309859154929577
$ws.Range("K11").Value = 0.1901408450704225
$ws.Range("L11").Value = 0.5633802816901409
$ws.Range("S11").Value = 0.01408450704225352
$ws.Range("G12").Value = 0.7962529274004684
$ws.Range("J12").Value = 0.1592505854800937
$ws.Range("K12").Value = 0.00468384074941452
$ws.Range("L12").Value = 0.01873536299765808
$ws.Range("S12").Value = 0.02107728337236534
$ws.Range("G13").Value = 0.728
$ws.Range("J13").Value = 0.216
$ws.Range("S13").Value = 0.056
$ws.Range("G14").Value = 0.875
$ws.Range("J14").Value = 0.125
$ws.Range("F15").Value = 0.02347417840375587
$ws.Range("H15").Value = 0.215962441314554
$ws.Range("I15").Value = 0.05868544600938967
$ws.Range("J15").Value = 0.3075117370892019
$ws.Range("K15").Value = 0.06338028169014084
$ws.Range("M15").Value = 0.02347417840375587
$ws.Range("N15").Value = 0.002347417840375587
$ws.Range("O15").Value = 0.07276995305164319
$ws.Range("S15").Value = 0.2323943661971831
$ws.Range("F16").Value = 0.02277432712215321
$ws.Range("H16").Value = 0.1780538302277433
$ws.Range("I16").Value = 0.08488612836438923
$ws.Range("J16").Value = 0.3478260869565217
$ws.Range("K16").Value = 0.1221532091097308
$ws.Range("M16").Value = 0.02898550724637681
$ws.Range("N16").Value = 0.006211180124223602
$ws.Range("O16").Value = 0.06418219461697723
$ws.Range("S16").Value = 0.1449275362318841
$ws.Range("F17").Value = 0.01935483870967742
$ws.Range("H17").Value = 0.2161290322580645
$ws.Range("I17").Value = 0.05591397849462366
$ws.Range("J17").Value = 0.4053763440860215
$ws.Range("K17").Value = 0.1053763440860215
$ws.Range("M17").Value = 0.02795698924731183
$ws.Range("N17").Value = 0.001075268817204301
$ws.Range("O17").Value = 0.05376344086021505
$ws.Range("S17").Value = 0.1150537634408602
$ws.Range("F18").Value = 0.01587301587301587
$ws.Range("H18").Value = 0.1927437641723356
$ws.Range("I18").Value = 0.07029478458049887
$ws.Range("J18").Value = 0.3968253968253968
$ws.Range("K18").Value = 0.1179138321995465
$ws.Range("M18").Value = 0.02494331065759637
$ws.Range("N18").Value = 0.002267573696145125
$ws.Range("O18").Value = 0.03854875283446712
$ws.Range("S18").Value = 0.1405895691609977
$ws.Range("F19").Value = 0.01151560178306092
$ws.Range("H19").Value = 0.2191679049034175
$ws.Range("I19").Value = 0.07652303120356613
$ws.Range("J19").Value = 0.3528974739970283
$ws.Range("K19").Value = 0.1210995542347697
$ws.Range("M19").Value = 0.02451708766716196
$ws.Range("N19").Value = 0.001485884101040119
$ws.Range("O19").Value = 0.06537890044576523
$ws.Range("S19").Value = 0.1274145616641902
